# threaded server working fine
#
# 1. Type a trailing space after "A 'purchase order' is characterised by:"
#    and move the (hidden) "_GoBack" last-edit bookmark there.
# 2. The old "_GoBack" location (inside the GUI-frameworks paragraph)
#    loses its bookmark and the two runs split around it collapse back
#    into a single run.
# 3. Refresh the cached "Last updated" DATE field result in the footer.

$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# ---------------------------------------------------------------------
# Step 1: insert the trailing space after "is characterised by:" (the
# occurrence that follows "'purchase order'").
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("'purchase order' is characterised by:") | Out-Null
$target.Collapse(0)
$target.InsertAfter(" ")

# $target now spans exactly the newly-inserted space character (Word
# COM semantics: InsertAfter grows the range to cover the inserted
# text). Use a throw-away marker placed right after it so we can get a
# reliable "end of space" insertion point for the bookmark -- adding a
# bookmark exactly at a run/paragraph boundary is ambiguous, so we
# anchor relative to real text instead and delete the marker after.
$target.Collapse(0)
$target.InsertAfter("@@GOBACKMARK@@")

$markerRange = $d.Content
$markerRange.Find.Execute("@@GOBACKMARK@@") | Out-Null
$markerRange.Collapse(1)

# Drop the old bookmark (if present) before adding the new one so the
# bookmark id / name is simply relocated.
$oldBookmarkExists = $false
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") { $oldBookmarkExists = $true }
}
if ($oldBookmarkExists) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null

$markerDel = $d.Content
$markerDel.Find.Execute("@@GOBACKMARK@@") | Out-Null
$markerDel.Delete()

# ---------------------------------------------------------------------
# Step 2: the GUI-frameworks paragraph used to hold the "_GoBack"
# bookmark between two runs; with it gone Word collapses the runs back
# into one contiguous run of text.
# ---------------------------------------------------------------------
$mergeLead = $d.Content
$mergeLead.Find.Execute(".  You are expected to ") | Out-Null
$mergeStart = $mergeLead.Start

$mergeTail = $d.Content
$mergeTail.Find.Execute("use one of the Java GUI frameworks for building your GUI; choose between AWT or Swing.") | Out-Null
$mergeEnd = $mergeTail.End

$mergedRange = $d.Range($mergeStart, $mergeEnd)
$mergedText = $mergedRange.Text
$mergedRange.Delete()

$reinsertPoint = $d.Range($mergeStart, $mergeStart)
$reinsertPoint.InsertAfter($mergedText)

# ---------------------------------------------------------------------
# Step 3: refresh the cached "Last updated" DATE field result shown in
# the footer.
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("03/03/2019 09:27", $true, $false, $false, $false, $false, $true, 1, $false, "04/03/2019 11:12", 2) | Out-Null
